$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four "Resolving-Mac" sending-cluster rows (old rows 14-17);
# Resolving-Mac no longer appears as a sending cluster after TPM recalculation.
$ws.Rows("14:17").Delete()

# Update recalculated TPM-based numeric values for the remaining rows (2-13)
# Row 2
$ws.Range("G2").Value = 1.008000333333333
$ws.Range("H2").Value = 3.024001
$ws.Range("I2").Value = 0.3525296793986107
$ws.Range("J2").Value = 0.3525296793986107
$ws.Range("M2").Value = 8.813278666666667
$ws.Range("N2").Value = 26.439836
$ws.Range("O2").Value = 0.3770976991891536
$ws.Range("P2").Value = 0.3770976991891536
$ws.Range("Q2").Value = 8.883787833759557
$ws.Range("R2").Value = 79.954090503836
$ws.Range("S2").Value = 0.1329381309971061
$ws.Range("T2").Value = 0.1329381309971061
# Row 3
$ws.Range("G3").Value = 1.008000333333333
$ws.Range("H3").Value = 3.024001
$ws.Range("I3").Value = 0.3525296793986107
$ws.Range("J3").Value = 0.3525296793986107
$ws.Range("O3").Value = 0.5522024902836482
$ws.Range("P3").Value = 0.5522024902836482
$ws.Range("Q3").Value = 13.00896233390411
$ws.Range("R3").Value = 117.080661005137
$ws.Range("S3").Value = 0.194667766862809
$ws.Range("T3").Value = 0.1946677668628089
# Row 4
$ws.Range("G4").Value = 1.008000333333333
$ws.Range("H4").Value = 3.024001
$ws.Range("I4").Value = 0.3525296793986107
$ws.Range("J4").Value = 0.3525296793986107
$ws.Range("M4").Value = 1.649921333333333
$ws.Range("N4").Value = 4.949764
$ws.Range("O4").Value = 0.07059592260441032
$ws.Range("P4").Value = 0.07059592260441033
$ws.Range("Q4").Value = 1.663121253973778
$ws.Range("R4").Value = 14.968091285764
$ws.Range("S4").Value = 0.02488715796258191
$ws.Range("T4").Value = 0.02488715796258191
# Row 5
$ws.Range("G5").Value = 1.008000333333333
$ws.Range("H5").Value = 3.024001
$ws.Range("I5").Value = 0.3525296793986107
$ws.Range("J5").Value = 0.3525296793986107
$ws.Range("M5").Value = 0.002428
$ws.Range("N5").Value = 0.007284
$ws.Range("O5").Value = 0.0001038879227879399
$ws.Range("P5").Value = 0.0001038879227879399
$ws.Range("Q5").Value = 0.002447424809333334
$ws.Range("R5").Value = 0.022026823284
$ws.Range("S5").Value = 0.0000366235761138201
$ws.Range("T5").Value = 0.00003662357611382009
# Row 6
$ws.Range("I6").Value = 0.5377259091975243
$ws.Range("J6").Value = 0.5377259091975243
$ws.Range("M6").Value = 8.813278666666667
$ws.Range("N6").Value = 26.439836
$ws.Range("O6").Value = 0.3770976991891536
$ws.Range("P6").Value = 0.3770976991891536
$ws.Range("Q6").Value = 13.55075379234889
$ws.Range("R6").Value = 121.95678413114
$ws.Range("S6").Value = 0.2027752031527822
$ws.Range("T6").Value = 0.2027752031527822
# Row 7
$ws.Range("I7").Value = 0.5377259091975243
$ws.Range("J7").Value = 0.5377259091975243
$ws.Range("O7").Value = 0.5522024902836482
$ws.Range("P7").Value = 0.5522024902836482
$ws.Range("S7").Value = 0.2969335861489118
$ws.Range("T7").Value = 0.2969335861489118
# Row 8
$ws.Range("I8").Value = 0.5377259091975243
$ws.Range("J8").Value = 0.5377259091975243
$ws.Range("M8").Value = 1.649921333333333
$ws.Range("N8").Value = 4.949764
$ws.Range("O8").Value = 0.07059592260441032
$ws.Range("P8").Value = 0.07059592260441033
$ws.Range("Q8").Value = 2.536817296984444
$ws.Range("R8").Value = 22.83135567286
$ws.Range("S8").Value = 0.03796125666809459
$ws.Range("T8").Value = 0.0379612566680946
# Row 9
$ws.Range("I9").Value = 0.5377259091975243
$ws.Range("J9").Value = 0.5377259091975243
$ws.Range("M9").Value = 0.002428
$ws.Range("N9").Value = 0.007284
$ws.Range("O9").Value = 0.0001038879227879399
$ws.Range("P9").Value = 0.0001038879227879399
$ws.Range("Q9").Value = 0.003733143073333333
$ws.Range("R9").Value = 0.03359828766
$ws.Range("S9").Value = 0.00005586322773578721
$ws.Range("T9").Value = 0.00005586322773578721
# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.313796
$ws.Range("H10").Value = 0.941388
$ws.Range("I10").Value = 0.1097444114038651
$ws.Range("J10").Value = 0.1097444114038651
$ws.Range("M10").Value = 8.813278666666667
$ws.Range("N10").Value = 26.439836
$ws.Range("O10").Value = 0.3770976991891536
$ws.Range("P10").Value = 0.3770976991891536
$ws.Range("Q10").Value = 2.765571592485334
$ws.Range("R10").Value = 24.890144332368
$ws.Range("S10").Value = 0.04138436503926543
$ws.Range("T10").Value = 0.04138436503926542
# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.313796
$ws.Range("H11").Value = 0.941388
$ws.Range("I11").Value = 0.1097444114038651
$ws.Range("J11").Value = 0.1097444114038651
$ws.Range("O11").Value = 0.5522024902836482
$ws.Range("P11").Value = 0.5522024902836482
$ws.Range("Q11").Value = 4.049760907350667
$ws.Range("R11").Value = 36.447848166156
$ws.Range("S11").Value = 0.06060113727192749
$ws.Range("T11").Value = 0.06060113727192748
# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.313796
$ws.Range("H12").Value = 0.941388
$ws.Range("I12").Value = 0.1097444114038651
$ws.Range("J12").Value = 0.1097444114038651
$ws.Range("M12").Value = 1.649921333333333
$ws.Range("N12").Value = 4.949764
$ws.Range("O12").Value = 0.07059592260441032
$ws.Range("P12").Value = 0.07059592260441033
$ws.Range("Q12").Value = 0.5177387147146667
$ws.Range("R12").Value = 4.659648432432
$ws.Range("S12").Value = 0.007747507973733823
$ws.Range("T12").Value = 0.007747507973733824
# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.313796
$ws.Range("H13").Value = 0.941388
$ws.Range("I13").Value = 0.1097444114038651
$ws.Range("J13").Value = 0.1097444114038651
$ws.Range("M13").Value = 0.002428
$ws.Range("N13").Value = 0.007284
$ws.Range("O13").Value = 0.0001038879227879399
$ws.Range("P13").Value = 0.0001038879227879399
$ws.Range("Q13").Value = 0.000761896688
$ws.Range("R13").Value = 0.006857070192
$ws.Range("S13").Value = 0.00001140111893833265
$ws.Range("T13").Value = 0.00001140111893833265
